$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing A:C data
# (the 0.1/1/0.5 ... progress-bar values) right into B:D.
$ws.Columns("A:A").Insert()

# Fill the freed-up column A with the row index (1-10) — the new
# "click to navigate" column used to jump to a given progress step.
for ($i = 1; $i -le 10; $i++) {
    $ws.Cells.Item($i, 1).Value = $i
}

# Match the author's recorded selection after making the edit.
$ws.Range("L14").Select()
